$wb = $excel.ActiveWorkbook

# ALC row 21
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 12500  # H21: 20000 -> 12500
$ws.Cells.Item(21, 9).Value = 12500  # I21: 20000 -> 12500
$ws.Cells.Item(21, 11).Value = 12500  # K21: 20000 -> 12500
$ws.Cells.Item(21, 13).Value = -12032  # M21: -19532 -> -12032

# ALC row 23
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(23, 8).Value = 12500  # H23: 20000 -> 12500
$ws.Cells.Item(23, 9).Value = 12500  # I23: 20000 -> 12500
$ws.Cells.Item(23, 11).Value = 12500  # K23: 20000 -> 12500
$ws.Cells.Item(23, 13).Value = -12266  # M23: -19766 -> -12266

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2750  # H40: 5266.1333 -> 2750
$ws.Cells.Item(40, 9).Value = 2666  # I40: 7249.125 -> 2666
$ws.Cells.Item(40, 10).Value = 2778  # J40: 2999.8572 -> 2778
$ws.Cells.Item(40, 11).Value = 2666  # K40: 7249.125 -> 2666
$ws.Cells.Item(40, 12).Value = 2778  # L40: 2999.8572 -> 2778
$ws.Cells.Item(40, 13).Value = -2491  # M40: -7074.125 -> -2491
$ws.Cells.Item(40, 14).Value = -3128  # N40: -3349.8572 -> -3128

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 11383.818  # H62: 11935.429 -> 11383.818
$ws.Cells.Item(62, 9).Value = 9750.1  # I62: 9889.6 -> 9750.1
$ws.Cells.Item(62, 10).Value = 12745.25  # J62: 13795.272 -> 12745.25
$ws.Cells.Item(62, 11).Value = 9750.1  # K62: 9889.6 -> 9750.1
$ws.Cells.Item(62, 12).Value = 12745.25  # L62: 13795.272 -> 12745.25
$ws.Cells.Item(62, 13).Value = -9126.1  # M62: -9265.6 -> -9126.1
$ws.Cells.Item(62, 14).Value = -13993.25  # N62: -15043.272 -> -13993.25

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 4848.5557  # H64: 5104.625 -> 4848.5557
$ws.Cells.Item(64, 9).Value = 4829.625  # I64: 5119.5713 -> 4829.625
$ws.Cells.Item(64, 11).Value = 4829.625  # K64: 5119.5713 -> 4829.625
$ws.Cells.Item(64, 13).Value = -4581.625  # M64: -4871.5713 -> -4581.625

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 11383.818  # H65: 11935.429 -> 11383.818
$ws.Cells.Item(65, 9).Value = 9750.1  # I65: 9889.6 -> 9750.1
$ws.Cells.Item(65, 10).Value = 12745.25  # J65: 13795.272 -> 12745.25
$ws.Cells.Item(65, 11).Value = 48750.5  # K65: 49448 -> 48750.5
$ws.Cells.Item(65, 12).Value = 63726.25  # L65: 68976.36 -> 63726.25
$ws.Cells.Item(65, 13).Value = -45630.5  # M65: -46328 -> -45630.5
$ws.Cells.Item(65, 14).Value = -69966.25  # N65: -75216.36 -> -69966.25

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 4848.5557  # H67: 5104.625 -> 4848.5557
$ws.Cells.Item(67, 9).Value = 4829.625  # I67: 5119.5713 -> 4829.625
$ws.Cells.Item(67, 11).Value = 4829.625  # K67: 5119.5713 -> 4829.625
$ws.Cells.Item(67, 13).Value = -3971.625  # M67: -4261.5713 -> -3971.625

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 6739.4614  # H74: 7184.357 -> 6739.4614
$ws.Cells.Item(74, 9).Value = 4749.25  # I74: 6193 -> 4749.25
$ws.Cells.Item(74, 10).Value = 7624  # J74: 7735.1113 -> 7624
$ws.Cells.Item(74, 11).Value = 4749.25  # K74: 6193 -> 4749.25
$ws.Cells.Item(74, 12).Value = 7624  # L74: 7735.1113 -> 7624
$ws.Cells.Item(74, 13).Value = -3813.25  # M74: -5257 -> -3813.25
$ws.Cells.Item(74, 14).Value = -9496  # N74: -9607.1113 -> -9496

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 6739.4614  # H77: 7184.357 -> 6739.4614
$ws.Cells.Item(77, 9).Value = 4749.25  # I77: 6193 -> 4749.25
$ws.Cells.Item(77, 10).Value = 7624  # J77: 7735.1113 -> 7624
$ws.Cells.Item(77, 11).Value = 23746.25  # K77: 30965 -> 23746.25
$ws.Cells.Item(77, 12).Value = 38120  # L77: 38675.5565 -> 38120
$ws.Cells.Item(77, 13).Value = -19066.25  # M77: -26285 -> -19066.25
$ws.Cells.Item(77, 14).Value = -47480  # N77: -48035.5565 -> -47480

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 363.23077  # H107: 545.3333 -> 363.23077
$ws.Cells.Item(107, 9).Value = 280.6  # I107: 371.7143 -> 280.6
$ws.Cells.Item(107, 10).Value = 638.6667  # J107: 1153 -> 638.6667
$ws.Cells.Item(107, 11).Value = 280.6  # K107: 371.7143 -> 280.6
$ws.Cells.Item(107, 12).Value = 638.6667  # L107: 1153 -> 638.6667
$ws.Cells.Item(107, 13).Value = 1639.4  # M107: 1548.2857 -> 1639.4
$ws.Cells.Item(107, 14).Value = -4478.6667  # N107: -4993 -> -4478.6667

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 22629.658  # H116: 24268.21 -> 22629.658
$ws.Cells.Item(116, 9).Value = 24486.605  # I116: 27497.758 -> 24486.605
$ws.Cells.Item(116, 10).Value = 14969.75  # J116: 13861.889 -> 14969.75
$ws.Cells.Item(116, 11).Value = 24486.605  # K116: 27497.758 -> 24486.605
$ws.Cells.Item(116, 12).Value = 14969.75  # L116: 13861.889 -> 14969.75
$ws.Cells.Item(116, 13).Value = -21044.605  # M116: -24055.758 -> -21044.605
$ws.Cells.Item(116, 14).Value = -21853.75  # N116: -20745.889 -> -21853.75

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(127, 8).Value = 5958.5  # H127: 3106.5557 -> 5958.5
$ws.Cells.Item(127, 9).Value = 7446.3335  # I127: 3567.1428 -> 7446.3335
$ws.Cells.Item(127, 10).Value = 1495  # J127: 1494.5 -> 1495
$ws.Cells.Item(127, 11).Value = 22339.0005  # K127: 10701.4284 -> 22339.0005
$ws.Cells.Item(127, 12).Value = 4485  # L127: 4483.5 -> 4485
$ws.Cells.Item(127, 13).Value = -17379.0005  # M127: -5741.428400000001 -> -17379.0005
$ws.Cells.Item(127, 14).Value = -14405  # N127: -14403.5 -> -14405

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 40275.44  # H132: 40940.95 -> 40275.44
$ws.Cells.Item(132, 9).Value = 51350.426  # I132: 52459.24 -> 51350.426
$ws.Cells.Item(132, 11).Value = 154051.278  # K132: 157377.72 -> 154051.278
$ws.Cells.Item(132, 13).Value = -151521.278  # M132: -154847.72 -> -151521.278

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1117797.2  # H137: 912147.8 -> 1117797.2
$ws.Cells.Item(137, 9).Value = 1878  # I137: 1628.3077 -> 1878
$ws.Cells.Item(137, 11).Value = 5634  # K137: 4884.9231 -> 5634
$ws.Cells.Item(137, 13).Value = -3084  # M137: -2334.9231 -> -3084

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5212623.5  # H32: 5267492.5 -> 5212623.5
$ws.Cells.Item(32, 9).Value = 5322328  # I32: 5379557 -> 5322328
$ws.Cells.Item(32, 11).Value = 5322328  # K32: 5379557 -> 5322328
$ws.Cells.Item(32, 13).Value = -5322041  # M32: -5379270 -> -5322041

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 5631.421  # H45: 5657.7896 -> 5631.421
$ws.Cells.Item(45, 9).Value = 2692.7144  # I45: 2728.5 -> 2692.7144
$ws.Cells.Item(45, 11).Value = 2692.7144  # K45: 2728.5 -> 2692.7144
$ws.Cells.Item(45, 13).Value = -2315.7144  # M45: -2351.5 -> -2315.7144

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 15781.947  # H61: 17285.766 -> 15781.947
$ws.Cells.Item(61, 9).Value = 15821.833  # I61: 18386.3 -> 15821.833
$ws.Cells.Item(61, 11).Value = 15821.833  # K61: 18386.3 -> 15821.833
$ws.Cells.Item(61, 13).Value = -15609.833  # M61: -18174.3 -> -15609.833

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 2459.475  # H74: 2455.5854 -> 2459.475
$ws.Cells.Item(74, 9).Value = 2118  # I74: 2105 -> 2118
$ws.Cells.Item(74, 10).Value = 3210.72  # J74: 3302.8333 -> 3210.72
$ws.Cells.Item(74, 11).Value = 2118  # K74: 2105 -> 2118
$ws.Cells.Item(74, 12).Value = 3210.72  # L74: 3302.8333 -> 3210.72
$ws.Cells.Item(74, 13).Value = -1244  # M74: -1231 -> -1244
$ws.Cells.Item(74, 14).Value = -4958.719999999999  # N74: -5050.8333 -> -4958.719999999999

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 2459.475  # H77: 2455.5854 -> 2459.475
$ws.Cells.Item(77, 9).Value = 2118  # I77: 2105 -> 2118
$ws.Cells.Item(77, 10).Value = 3210.72  # J77: 3302.8333 -> 3210.72
$ws.Cells.Item(77, 11).Value = 10590  # K77: 10525 -> 10590
$ws.Cells.Item(77, 12).Value = 16053.6  # L77: 16514.1665 -> 16053.6
$ws.Cells.Item(77, 13).Value = -6222  # M77: -6157 -> -6222
$ws.Cells.Item(77, 14).Value = -24789.6  # N77: -25250.1665 -> -24789.6

# ARM row 92
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(92, 8).Value = 70439.8  # H92: 69199.664 -> 70439.8
$ws.Cells.Item(92, 10).Value = 70439.8  # J92: 69199.664 -> 70439.8
$ws.Cells.Item(92, 12).Value = 70439.8  # L92: 69199.664 -> 70439.8
$ws.Cells.Item(92, 14).Value = -75431.8  # N92: -74191.664 -> -75431.8

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2801.182  # H122: 2991.4 -> 2801.182
$ws.Cells.Item(122, 9).Value = 2351.625  # I122: 2559.1428 -> 2351.625
$ws.Cells.Item(122, 11).Value = 7054.875  # K122: 7677.428400000001 -> 7054.875
$ws.Cells.Item(122, 13).Value = -4604.875  # M122: -5227.428400000001 -> -4604.875

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 15781.947  # H136: 17285.766 -> 15781.947
$ws.Cells.Item(136, 9).Value = 15821.833  # I136: 18386.3 -> 15821.833
$ws.Cells.Item(136, 11).Value = 47465.499  # K136: 55158.89999999999 -> 47465.499
$ws.Cells.Item(136, 13).Value = -44915.499  # M136: -52608.89999999999 -> -44915.499

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 2464.9167  # H64: 3509.7273 -> 2464.9167
$ws.Cells.Item(64, 9).Value = 350  # I64: 305.25 -> 350
$ws.Cells.Item(64, 10).Value = 4579.8335  # J64: 5340.857 -> 4579.8335
$ws.Cells.Item(64, 11).Value = 350  # K64: 305.25 -> 350
$ws.Cells.Item(64, 12).Value = 4579.8335  # L64: 5340.857 -> 4579.8335
$ws.Cells.Item(64, 13).Value = -125  # M64: -80.25 -> -125
$ws.Cells.Item(64, 14).Value = -5029.8335  # N64: -5790.857 -> -5029.8335

# BSM row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(67, 8).Value = 2464.9167  # H67: 3509.7273 -> 2464.9167
$ws.Cells.Item(67, 9).Value = 350  # I67: 305.25 -> 350
$ws.Cells.Item(67, 10).Value = 4579.8335  # J67: 5340.857 -> 4579.8335
$ws.Cells.Item(67, 11).Value = 350  # K67: 305.25 -> 350
$ws.Cells.Item(67, 12).Value = 4579.8335  # L67: 5340.857 -> 4579.8335
$ws.Cells.Item(67, 13).Value = 430  # M67: 474.75 -> 430
$ws.Cells.Item(67, 14).Value = -6139.8335  # N67: -6900.857 -> -6139.8335

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5719.146  # H31: 5375.577 -> 5719.146
$ws.Cells.Item(31, 9).Value = 2055.56  # I31: 1942.9286 -> 2055.56
$ws.Cells.Item(31, 10).Value = 9701.305  # J31: 9380.333000000001 -> 9701.305
$ws.Cells.Item(31, 11).Value = 2055.56  # K31: 1942.9286 -> 2055.56
$ws.Cells.Item(31, 12).Value = 9701.305  # L31: 9380.333000000001 -> 9701.305
$ws.Cells.Item(31, 13).Value = -1760.56  # M31: -1647.9286 -> -1760.56
$ws.Cells.Item(31, 14).Value = -10291.305  # N31: -9970.333000000001 -> -10291.305

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 5719.146  # H34: 5375.577 -> 5719.146
$ws.Cells.Item(34, 9).Value = 2055.56  # I34: 1942.9286 -> 2055.56
$ws.Cells.Item(34, 10).Value = 9701.305  # J34: 9380.333000000001 -> 9701.305
$ws.Cells.Item(34, 11).Value = 2055.56  # K34: 1942.9286 -> 2055.56
$ws.Cells.Item(34, 12).Value = 9701.305  # L34: 9380.333000000001 -> 9701.305
$ws.Cells.Item(34, 13).Value = -1853.56  # M34: -1740.9286 -> -1853.56
$ws.Cells.Item(34, 14).Value = -10105.305  # N34: -9784.333000000001 -> -10105.305

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2354.4  # H58: 1746.8462 -> 2354.4
$ws.Cells.Item(58, 9).Value = 1483.8182  # I58: 1519 -> 1483.8182
$ws.Cells.Item(58, 10).Value = 4748.5  # J58: 3000 -> 4748.5
$ws.Cells.Item(58, 11).Value = 1483.8182  # K58: 1519 -> 1483.8182
$ws.Cells.Item(58, 12).Value = 4748.5  # L58: 3000 -> 4748.5
$ws.Cells.Item(58, 13).Value = -1280.8182  # M58: -1316 -> -1280.8182
$ws.Cells.Item(58, 14).Value = -5154.5  # N58: -3406 -> -5154.5

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 6840.6924  # H99: 6230.9414 -> 6840.6924
$ws.Cells.Item(99, 9).Value = 6018.4287  # I99: 5375.091 -> 6018.4287
$ws.Cells.Item(99, 11).Value = 6018.4287  # K99: 5375.091 -> 6018.4287
$ws.Cells.Item(99, 13).Value = -4520.4287  # M99: -3877.091 -> -4520.4287

# CRP row 100
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(100, 8).Value = 68916.25  # H100: 70449.75 -> 68916.25
$ws.Cells.Item(100, 10).Value = 68916.25  # J100: 70449.75 -> 68916.25
$ws.Cells.Item(100, 12).Value = 68916.25  # L100: 70449.75 -> 68916.25
$ws.Cells.Item(100, 14).Value = -71080.25  # N100: -72613.75 -> -71080.25

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 6840.6924  # H126: 6230.9414 -> 6840.6924
$ws.Cells.Item(126, 9).Value = 6018.4287  # I126: 5375.091 -> 6018.4287
$ws.Cells.Item(126, 11).Value = 18055.2861  # K126: 16125.273 -> 18055.2861
$ws.Cells.Item(126, 13).Value = -15585.2861  # M126: -13655.273 -> -15585.2861

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 2354.4  # H136: 1746.8462 -> 2354.4
$ws.Cells.Item(136, 9).Value = 1483.8182  # I136: 1519 -> 1483.8182
$ws.Cells.Item(136, 10).Value = 4748.5  # J136: 3000 -> 4748.5
$ws.Cells.Item(136, 11).Value = 4451.4546  # K136: 4557 -> 4451.4546
$ws.Cells.Item(136, 12).Value = 14245.5  # L136: 9000 -> 14245.5
$ws.Cells.Item(136, 13).Value = -1901.4546  # M136: -2007 -> -1901.4546
$ws.Cells.Item(136, 14).Value = -19345.5  # N136: -14100 -> -19345.5

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 277813.34  # H2: 52894776 -> 277813.34
$ws.Cells.Item(2, 9).Value = 384645.94  # I2: 416698.16 -> 384645.94
$ws.Cells.Item(2, 10).Value = 48.6  # J2: 142857180 -> 48.6
$ws.Cells.Item(2, 11).Value = 2307875.64  # K2: 2500188.96 -> 2307875.64
$ws.Cells.Item(2, 12).Value = 291.6  # L2: 857143080 -> 291.6
$ws.Cells.Item(2, 13).Value = -2307762.64  # M2: -2500075.96 -> -2307762.64
$ws.Cells.Item(2, 14).Value = -517.6  # N2: -857143306 -> -517.6

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 43.5  # H38: 41.52941 -> 43.5
$ws.Cells.Item(38, 9).Value = 31.222221  # I38: 29.1 -> 31.222221
$ws.Cells.Item(38, 11).Value = 93.666663  # K38: 87.30000000000001 -> 93.666663
$ws.Cells.Item(38, 13).Value = 253.333337  # M38: 259.7 -> 253.333337

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 4666.3335  # H39: 3979.8 -> 4666.3335
$ws.Cells.Item(39, 9).Value = 2500  # I39: 2750 -> 2500
$ws.Cells.Item(39, 10).Value = 5749.5  # J39: 4799.6665 -> 5749.5
$ws.Cells.Item(39, 11).Value = 7500  # K39: 8250 -> 7500
$ws.Cells.Item(39, 12).Value = 17248.5  # L39: 14398.9995 -> 17248.5
$ws.Cells.Item(39, 13).Value = -7206  # M39: -7956 -> -7206
$ws.Cells.Item(39, 14).Value = -17836.5  # N39: -14986.9995 -> -17836.5

# CUL row 63
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 13177.333  # H63: 10014 -> 13177.333
$ws.Cells.Item(63, 10).Value = 13177.333  # J63: 10014 -> 13177.333
$ws.Cells.Item(63, 12).Value = 39531.999  # L63: 30042 -> 39531.999
$ws.Cells.Item(63, 14).Value = -41029.999  # N63: -31540 -> -41029.999

# CUL row 66
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(66, 8).Value = 13177.333  # H66: 10014 -> 13177.333
$ws.Cells.Item(66, 10).Value = 13177.333  # J66: 10014 -> 13177.333
$ws.Cells.Item(66, 12).Value = 118595.997  # L66: 90126 -> 118595.997
$ws.Cells.Item(66, 14).Value = -126083.997  # N66: -97614 -> -126083.997

# CUL row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(114, 8).Value = 1850.25  # H114: 2016.0834 -> 1850.25
$ws.Cells.Item(114, 9).Value = 266.07693  # I114: 572.2308 -> 266.07693
$ws.Cells.Item(114, 11).Value = 798.2307900000001  # K114: 1716.6924 -> 798.2307900000001
$ws.Cells.Item(114, 13).Value = 2455.76921  # M114: 1537.3076 -> 2455.76921

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 2000  # H139: 1165.6 -> 2000
$ws.Cells.Item(139, 9).Value = 2000  # I139: 1165.6 -> 2000
$ws.Cells.Item(139, 11).Value = 6000  # K139: 3496.8 -> 6000
$ws.Cells.Item(139, 13).Value = -860  # M139: 1643.2 -> -860

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 6719.4375  # H122: 6655.0586 -> 6719.4375
$ws.Cells.Item(122, 9).Value = 7775.9165  # I122: 7953.0835 -> 7775.9165
$ws.Cells.Item(122, 10).Value = 3550  # J122: 3539.8 -> 3550
$ws.Cells.Item(122, 11).Value = 23327.7495  # K122: 23859.2505 -> 23327.7495
$ws.Cells.Item(122, 12).Value = 10650  # L122: 10619.4 -> 10650
$ws.Cells.Item(122, 13).Value = -20877.7495  # M122: -21409.2505 -> -20877.7495
$ws.Cells.Item(122, 14).Value = -15550  # N122: -15519.4 -> -15550

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2989.2222  # H16: 3414.7144 -> 2989.2222
$ws.Cells.Item(16, 9).Value = 1750.25  # I16: 2000.3334 -> 1750.25
$ws.Cells.Item(16, 10).Value = 3980.4  # J16: 4475.5 -> 3980.4
$ws.Cells.Item(16, 11).Value = 1750.25  # K16: 2000.3334 -> 1750.25
$ws.Cells.Item(16, 12).Value = 3980.4  # L16: 4475.5 -> 3980.4
$ws.Cells.Item(16, 13).Value = -1580.25  # M16: -1830.3334 -> -1580.25
$ws.Cells.Item(16, 14).Value = -4320.4  # N16: -4815.5 -> -4320.4

# LTW row 57
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(57, 8).Value = 0  # H57: 44000 -> 0
$ws.Cells.Item(57, 9).Value = 0  # I57: 44000 -> 0
$ws.Cells.Item(57, 11).Value = 0  # K57: 44000 -> 0
$ws.Cells.Item(57, 13).ClearContents()  # M57: was -43434, removed

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2625.25  # H61: 2029.5385 -> 2625.25
$ws.Cells.Item(61, 9).Value = 3004  # I61: 1888.7 -> 3004
$ws.Cells.Item(61, 11).Value = 3004  # K61: 1888.7 -> 3004
$ws.Cells.Item(61, 13).Value = -2802  # M61: -1686.7 -> -2802

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 2625.25  # H113: 2029.5385 -> 2625.25
$ws.Cells.Item(113, 9).Value = 3004  # I113: 1888.7 -> 3004
$ws.Cells.Item(113, 11).Value = 3004  # K113: 1888.7 -> 3004
$ws.Cells.Item(113, 13).Value = -834  # M113: 281.3 -> -834

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 770978.0600000001  # H132: 937157.0600000001 -> 770978.0600000001
$ws.Cells.Item(132, 9).Value = 1019091.6  # I132: 1331919.8 -> 1019091.6
$ws.Cells.Item(132, 11).Value = 3057274.8  # K132: 3995759.4 -> 3057274.8
$ws.Cells.Item(132, 13).Value = -3054744.8  # M132: -3993229.4 -> -3054744.8

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 5406.593  # H136: 5369.222 -> 5406.593
$ws.Cells.Item(136, 9).Value = 5003.6313  # I136: 4803.2 -> 5003.6313
$ws.Cells.Item(136, 10).Value = 6363.625  # J136: 6986.4287 -> 6363.625
$ws.Cells.Item(136, 11).Value = 15010.8939  # K136: 14409.6 -> 15010.8939
$ws.Cells.Item(136, 12).Value = 19090.875  # L136: 20959.2861 -> 19090.875
$ws.Cells.Item(136, 13).Value = -12460.8939  # M136: -11859.6 -> -12460.8939
$ws.Cells.Item(136, 14).Value = -24190.875  # N136: -26059.2861 -> -24190.875

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 4486.048  # H107: 3172.3215 -> 4486.048
$ws.Cells.Item(107, 9).Value = 1177.5  # I107: 810.26666 -> 1177.5
$ws.Cells.Item(107, 10).Value = 7493.8184  # J107: 5897.769 -> 7493.8184
$ws.Cells.Item(107, 11).Value = 3532.5  # K107: 2430.79998 -> 3532.5
$ws.Cells.Item(107, 12).Value = 22481.4552  # L107: 17693.307 -> 22481.4552
$ws.Cells.Item(107, 13).Value = -1612.5  # M107: -510.7999799999998 -> -1612.5
$ws.Cells.Item(107, 14).Value = -26321.4552  # N107: -21533.307 -> -26321.4552

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 1070  # H113: 809.76 -> 1070
$ws.Cells.Item(113, 9).Value = 772.4167  # I113: 638.0714 -> 772.4167
$ws.Cells.Item(113, 10).Value = 1427.1  # J113: 1028.2727 -> 1427.1
$ws.Cells.Item(113, 11).Value = 2317.2501  # K113: 1914.2142 -> 2317.2501
$ws.Cells.Item(113, 12).Value = 4281.299999999999  # L113: 3084.8181 -> 4281.299999999999
$ws.Cells.Item(113, 13).Value = -147.2501000000002  # M113: 255.7857999999999 -> -147.2501000000002
$ws.Cells.Item(113, 14).Value = -8621.299999999999  # N113: -7424.8181 -> -8621.299999999999

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2591.5264  # H126: 2654.9443 -> 2591.5264
$ws.Cells.Item(126, 9).Value = 2349.818  # I126: 2439.8 -> 2349.818
$ws.Cells.Item(126, 11).Value = 7049.454000000001  # K126: 7319.400000000001 -> 7049.454000000001
$ws.Cells.Item(126, 13).Value = -4579.454000000001  # M126: -4849.400000000001 -> -4579.454000000001

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 13577154  # H136: 15839922 -> 13577154
$ws.Cells.Item(136, 9).Value = 15839159  # I136: 17279028 -> 15839159
$ws.Cells.Item(136, 10).Value = 5127  # J136: 9749.5 -> 5127
$ws.Cells.Item(136, 11).Value = 47517477  # K136: 51837084 -> 47517477
$ws.Cells.Item(136, 12).Value = 15381  # L136: 29248.5 -> 15381
$ws.Cells.Item(136, 13).Value = -47514927  # M136: -51834534 -> -47514927
$ws.Cells.Item(136, 14).Value = -20481  # N136: -34348.5 -> -20481
